$wb = $excel.ActiveWorkbook

# "想去人数" (number of interested attendees) increased for two events that
# appear in both the "展览" and "全部类型" sheets.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1239
    $ws.Range("F6").Value = 161
}
